$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

$timeValues = @(
    "2021-10-05 13:39:29.278560",
    "2021-10-05 13:39:29.278573",
    "2021-10-05 13:39:29.278578",
    "2021-10-05 13:39:29.278581",
    "2021-10-05 13:39:29.278584",
    "2021-10-05 13:39:29.278587",
    "2021-10-05 13:39:29.278590",
    "2021-10-05 13:39:29.278593",
    "2021-10-05 13:39:29.278597",
    "2021-10-05 13:39:29.278600",
    "2021-10-05 13:39:29.278603",
    "2021-10-05 13:39:29.278606",
    "2021-10-05 13:39:29.278609",
    "2021-10-05 13:39:29.278611",
    "2021-10-05 13:39:29.278614",
    "2021-10-05 13:39:29.278617",
    "2021-10-05 13:39:29.278621",
    "2021-10-05 13:39:29.278624",
    "2021-10-05 13:39:29.278627",
    "2021-10-05 13:39:29.278630",
    "2021-10-05 13:39:29.278632",
    "2021-10-05 13:39:29.278635",
    "2021-10-05 13:39:29.278638",
    "2021-10-05 13:39:29.278641",
    "2021-10-05 13:39:29.278645",
    "2021-10-05 13:39:29.278648",
    "2021-10-05 13:39:29.278651",
    "2021-10-05 13:39:29.278654",
    "2021-10-05 13:39:29.278656",
    "2021-10-05 13:39:29.278659",
    "2021-10-05 13:39:29.278662",
    "2021-10-05 13:39:29.278666",
    "2021-10-05 13:39:29.278669",
    "2021-10-05 13:39:29.278672",
    "2021-10-05 13:39:29.278675",
    "2021-10-05 13:39:29.278678",
    "2021-10-05 13:39:29.278681",
    "2021-10-05 13:39:29.278684",
    "2021-10-05 13:39:29.278687",
    "2021-10-05 13:39:29.278690",
    "2021-10-05 13:39:29.278693",
    "2021-10-05 13:39:29.278696",
    "2021-10-05 13:39:29.278699",
    "2021-10-05 13:39:29.278702",
    "2021-10-05 13:39:29.278705",
    "2021-10-05 13:39:29.278708",
    "2021-10-05 13:39:29.278711",
    "2021-10-05 13:39:29.278714",
    "2021-10-05 13:39:29.278717",
    "2021-10-05 13:39:29.278720",
    "2021-10-05 13:39:29.278723",
    "2021-10-05 13:39:29.278726",
    "2021-10-05 13:39:29.278729",
    "2021-10-05 13:39:29.278732",
    "2021-10-05 13:39:29.278735",
    "2021-10-05 13:39:29.278738",
    "2021-10-05 13:39:29.278741",
    "2021-10-05 13:39:29.278743",
    "2021-10-05 13:39:29.278746",
    "2021-10-05 13:39:29.278749",
    "2021-10-05 13:39:29.278752",
    "2021-10-05 13:39:29.278755",
    "2021-10-05 13:39:29.278758",
    "2021-10-05 13:39:29.278761",
    "2021-10-05 13:39:29.278765",
    "2021-10-05 13:39:29.278769",
    "2021-10-05 13:39:29.278772",
    "2021-10-05 13:39:29.278775",
    "2021-10-05 13:39:29.278778",
    "2021-10-05 13:39:29.278780",
    "2021-10-05 13:39:29.278783",
    "2021-10-05 13:39:29.278786",
    "2021-10-05 13:39:29.278789",
    "2021-10-05 13:39:29.278792",
    "2021-10-05 13:39:29.278795",
    "2021-10-05 13:39:29.278798",
    "2021-10-05 13:39:29.278803",
    "2021-10-05 13:39:29.278807",
    "2021-10-05 13:39:29.278810",
    "2021-10-05 13:39:29.278813",
    "2021-10-05 13:39:29.278816",
    "2021-10-05 13:39:29.278819",
    "2021-10-05 13:39:29.278822",
    "2021-10-05 13:39:29.278825",
    "2021-10-05 13:39:29.278828"
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeValues[$i]
}
